$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data for rows 2..9 (columns A..T)
$data = @(
    @("ECs","Cck","Cckbr","FAPs",3,1,1.862015666666667,5.586047,0.399282427039184,0.3992824270391842,2,0.6666666666666666,0.08226166666666666,0.246785,0.9036565896241614,0.9036565896241615,0.1531725120994444,1.378552608895,0.3608141963150871,0.3608141963150873),
    @("ECs","Cck","Cckbr","MuSCs",3,1,1.862015666666667,5.586047,0.399282427039184,0.3992824270391842,1,0.3333333333333333,0.008770333333333333,0.026311,0.09634341037583853,0.09634341037583853,0.01633049806855555,0.146974482617,0.03846823072409691,0.03846823072409693),
    @("FAPs","Cck","Cckbr","FAPs",3,1,2.017090333333333,6.051271,0.4325359545939786,0.4325359545939787,2,0.6666666666666666,0.08226166666666666,0.246785,0.9036565896241614,0.9036565896241615,0.1659292126372222,1.493362913735,0.3908639656182259,0.390863965618226),
    @("FAPs","Cck","Cckbr","MuSCs",3,1,2.017090333333333,6.051271,0.4325359545939786,0.4325359545939787,1,0.3333333333333333,0.008770333333333333,0.026311,0.09634341037583853,0.09634341037583853,0.01769055458677778,0.159214991281,0.04167198897575274,0.04167198897575275),
    @("Inflammatory-Mac","Cck","Cckbr","FAPs",1,0.3333333333333333,0.1819486666666667,0.545846,0.0390162695855639,0.03901626958556391,2,0.6666666666666666,0.08226166666666666,0.246785,0.9036565896241614,0.9036565896241615,0.01496740056777778,0.13470660511,0.03525730911354757,0.03525730911354758),
    @("Inflammatory-Mac","Cck","Cckbr","MuSCs",1,0.3333333333333333,0.1819486666666667,0.545846,0.0390162695855639,0.03901626958556391,1,0.3333333333333333,0.008770333333333333,0.026311,0.09634341037583853,0.09634341037583853,0.001595750456222222,0.014361754106,0.003758960472016331,0.003758960472016331),
    @("MuSCs","Cck","Cckbr","FAPs",3,1,0.6023503333333334,1.807051,0.1291653487812732,0.1291653487812732,2,0.6666666666666666,0.08226166666666666,0.246785,0.9036565896241614,0.9036565896241615,0.04955034233722223,0.445953081035,0.1167211185773007,0.1167211185773007),
    @("MuSCs","Cck","Cckbr","MuSCs",3,1,0.6023503333333334,1.807051,0.1291653487812732,0.1291653487812732,1,0.3333333333333333,0.008770333333333333,0.026311,0.09634341037583853,0.09634341037583853,0.005282813206777779,0.047545318861,0.01244423020397252,0.01244423020397252)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# Row 10 no longer exists in the updated data - delete it so the sheet
# shrinks from 9 data rows + header (A1:T10) to 8 data rows + header (A1:T9)
$ws.Rows("10:10").Delete()

